# "Fixed tables on Resources page"
#
# The Human-sheet Reference column contained <a href="..."> links without
# target="_blank", so clicking them navigated away from the site instead of
# opening a new tab. Bring them in line with the other sheets' convention
# of opening external references in a new tab, and leave the workbook with
# the Human tab (the first / primary table) selected and active, with the
# last-edited cell (A9) as the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Human")

$ws.Range("A2").Value2 = '<a href="https://www.nature.com/articles/nn.3980" target="_blank">Johnson</a>'
$ws.Range("A3").Value2 = '<a href="http://www.pnas.org/content/112/23/7285" target="_blank">Darmanis</a>'
$ws.Range("A4").Value2 = '<a href="http://science.sciencemag.org/content/352/6293/1586.long"  target="_blank">Lake</a>'
$ws.Range("A5").Value2 = '<a href="https://doi.org/10.1016/j.cell.2015.09.004"  target="_blank">Pollen</a>'
$ws.Range("A6").Value2 = '<a href="http://science.sciencemag.org/content/358/6368/1318"  target="_blank">Nowakowski</a>'
$ws.Range("A7").Value2 = '<a href="https://www.nature.com/articles/nature25980"  target="_blank">Zhong</a>'
$ws.Range("A8").Value2 = '<a href="https://www.nature.com/articles/nbt.4038"  target="_blank">Lake</a>'
$ws.Range("A9").Value2 = '<a href="https://www.nature.com/articles/s41422-018-0053-3"  target="_blank">Fan</a>'

# Make the Human sheet the active tab (it was "Human organoid" before) and
# leave the cursor on the last cell that was touched.
$ws.Select()
$ws.Range("A9").Select()
